# B6-PowerPoint.pptx edit
#
# 1) Re-style the three tables (slides 14-16) from the Google-Slides-exported
#    custom "Table_0" style to the built-in PowerPoint table style
#    {D891B640-AC34-49DF-9732-228236058C9B}.
# 2) Re-colour the deck's (single) theme from the "Integral" / "Red Violet"
#    palette over to the standard Office "Office Theme" palette - this is the
#    part of the edit that is reachable through the exposed theme-colour
#    object model (ThemeColorScheme), matching the colour values the target
#    theme1.xml ends up with.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Table styles
# ---------------------------------------------------------------------
$builtinTableStyle = "{D891B640-AC34-49DF-9732-228236058C9B}"

for ($slideIdx = 1; $slideIdx -le $p.Slides.Count; $slideIdx++) {
    $slide = $p.Slides.Item($slideIdx)
    for ($shapeIdx = 1; $shapeIdx -le $slide.Shapes.Count; $shapeIdx++) {
        $shape = $slide.Shapes.Item($shapeIdx)
        if ($shape.HasTable) {
            $table = $shape.Table
            if ($table.Style.Name -eq "{694FB601-E31B-4E7A-99D3-87AB263675D6}") {
                $table.ApplyStyle($builtinTableStyle)
            }
        }
    }
}

# ---------------------------------------------------------------------
# 2. Theme colours: Integral / Red Violet -> Office Theme
# ---------------------------------------------------------------------
function ConvertTo-ComRGB($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Order matches ThemeColorScheme indices 1..12.
$officeThemeColors = @(
    "000000",  # dk1
    "FFFFFF",  # lt1
    "44546A",  # dk2
    "E7E6E6",  # lt2
    "5B9BD5",  # accent1
    "ED7D31",  # accent2
    "A5A5A5",  # accent3
    "FFC000",  # accent4
    "4472C4",  # accent5
    "70AD47",  # accent6
    "0563C1",  # hlink
    "954F72"   # folHlink
)

$themeColors = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Item($i).RGB = ConvertTo-ComRGB $officeThemeColors[$i - 1]
}
